$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 660 (pushes existing rows 660-701 down to 661-702)
$ws.Rows.Item(660).Insert()

# Populate the newly inserted row with the new data point for 2026/01/18.
# Force text formatting on column A first so the date-like string "2026/01/18"
# is stored as literal text (matching the rest of the column) instead of being
# auto-converted into a date serial number by Excel's type inference.
$ws.Range("A660").NumberFormat = "@"
$ws.Range("A660").Value = "2026/01/18"
$ws.Range("A660").Style = "Normal"

$ws.Range("B660").Value = "日"
$ws.Range("C660").Value = 4
$ws.Range("D660").Value = 201
